$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.247.10'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.840.61'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6700'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07423'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2963'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07725'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.025'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6787'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.753.28'
$ws.Range('E14').Value = '  -4.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.188'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '29.340.61'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008264'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9997'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.249'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.707'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1410'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.513'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.199'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.087'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.191'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05354'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.879'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7575'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.25%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '1.334.63'
$ws.Range('E37').Value = '  +3.51%  '
$ws.Range('E38').Value = '  -2.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.733'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9232'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.967'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.08005'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.85%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '2.012.22'
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000124'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5164'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.766'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '63.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.244'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.54%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05933'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.15%  '
